$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd34"
$ws.Range("C2").Value = "Sell"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 141.3574043333333
$ws.Range("H2").Value = 424.072213
$ws.Range("I2").Value = 0.4954750229273862
$ws.Range("J2").Value = 0.4954750229273862
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 74.89537566666667
$ws.Range("N2").Value = 224.686127
$ws.Range("O2").Value = 0.6660623326691122
$ws.Range("P2").Value = 0.6660623326691123
$ws.Range("Q2").Value = 10587.01590080989
$ws.Range("R2").Value = 95283.14310728904
$ws.Range("S2").Value = 0.3300172495502967
$ws.Range("T2").Value = 0.3300172495502967

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd34"
$ws.Range("C3").Value = "Sell"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 141.3574043333333
$ws.Range("H3").Value = 424.072213
$ws.Range("I3").Value = 0.4954750229273862
$ws.Range("J3").Value = 0.4954750229273862
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 37.54961933333333
$ws.Range("N3").Value = 112.648858
$ws.Range("O3").Value = 0.3339376673308877
$ws.Range("P3").Value = 0.3339376673308878
$ws.Range("Q3").Value = 5307.916722664751
$ws.Range("R3").Value = 47771.25050398274
$ws.Range("S3").Value = 0.1654577733770894
$ws.Range("T3").Value = 0.1654577733770895

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd34"
$ws.Range("C4").Value = "Sell"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 140.4344916666667
$ws.Range("H4").Value = 421.303475
$ws.Range("I4").Value = 0.4922401009448182
$ws.Range("J4").Value = 0.4922401009448182
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 74.89537566666667
$ws.Range("N4").Value = 224.686127
$ws.Range("O4").Value = 0.6660623326691122
$ws.Range("P4").Value = 0.6660623326691123
$ws.Range("Q4").Value = 10517.89400993237
$ws.Range("R4").Value = 94661.04608939133
$ws.Range("S4").Value = 0.3278625898685849
$ws.Range("T4").Value = 0.327862589868585

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd34"
$ws.Range("C5").Value = "Sell"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 140.4344916666667
$ws.Range("H5").Value = 421.303475
$ws.Range("I5").Value = 0.4922401009448182
$ws.Range("J5").Value = 0.4922401009448182
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 37.54961933333333
$ws.Range("N5").Value = 112.648858
$ws.Range("O5").Value = 0.3339376673308877
$ws.Range("P5").Value = 0.3339376673308878
$ws.Range("Q5").Value = 5273.261703353505
$ws.Range("R5").Value = 47459.35533018155
$ws.Range("S5").Value = 0.1643775110762333
$ws.Range("T5").Value = 0.1643775110762333

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cd34"
$ws.Range("C6").Value = "Sell"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.248853
$ws.Range("H6").Value = 0.7465590000000001
$ws.Range("I6").Value = 0.0008722602573388757
$ws.Range("J6").Value = 0.0008722602573388757
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 74.89537566666667
$ws.Range("N6").Value = 224.686127
$ws.Range("O6").Value = 0.6660623326691122
$ws.Range("P6").Value = 0.6660623326691123
$ws.Range("Q6").Value = 18.637938920777
$ws.Range("R6").Value = 167.741450286993
$ws.Range("S6").Value = 0.0005809797016976916
$ws.Range("T6").Value = 0.0005809797016976917

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd34"
$ws.Range("C7").Value = "Sell"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.248853
$ws.Range("H7").Value = 0.7465590000000001
$ws.Range("I7").Value = 0.0008722602573388757
$ws.Range("J7").Value = 0.0008722602573388757
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 37.54961933333333
$ws.Range("N7").Value = 112.648858
$ws.Range("O7").Value = 0.3339376673308877
$ws.Range("P7").Value = 0.3339376673308878
$ws.Range("Q7").Value = 9.344335419958
$ws.Range("R7").Value = 84.099018779622
$ws.Range("S7").Value = 0.000291280555641184
$ws.Range("T7").Value = 0.000291280555641184

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd34"
$ws.Range("C8").Value = "Sell"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.255981999999999
$ws.Range("H8").Value = 9.767945999999998
$ws.Range("I8").Value = 0.01141261587045664
$ws.Range("J8").Value = 0.01141261587045664
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 74.89537566666667
$ws.Range("N8").Value = 224.686127
$ws.Range("O8").Value = 0.6660623326691122
$ws.Range("P8").Value = 0.6660623326691123
$ws.Range("Q8").Value = 243.8579950539046
$ws.Range("R8").Value = 2194.721955485142
$ws.Range("S8").Value = 0.007601513548532881
$ws.Range("T8").Value = 0.007601513548532882

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd34"
$ws.Range("C9").Value = "Sell"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.255981999999999
$ws.Range("H9").Value = 9.767945999999998
$ws.Range("I9").Value = 0.01141261587045664
$ws.Range("J9").Value = 0.01141261587045664
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 37.54961933333333
$ws.Range("N9").Value = 112.648858
$ws.Range("O9").Value = 0.3339376673308877
$ws.Range("P9").Value = 0.3339376673308878
$ws.Range("Q9").Value = 122.2608846561853
$ws.Range("R9").Value = 1100.347961905668
$ws.Range("S9").Value = 0.00381110232192376
$ws.Range("T9").Value = 0.00381110232192376
